# OLX Monitor run @ 2026-02-15 22:01
# Updates the "last checked" timestamp on the summary sheet and appends a
# fresh pair of log rows (mirroring the existing 21:58 / 21:59 rows) to each
# per-profile detail sheet.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2026-02-15 22:01"

# --- 1) PODSUMOWANIE: bump "Data ostatniego sprawdzenia" for every profile ---
$summary = $wb.Worksheets.Item("PODSUMOWANIE")
$summary.Range("B2").Value = $newTimestamp
$summary.Range("B3").Value = $newTimestamp
$summary.Range("B4").Value = $newTimestamp
$summary.Range("B5").Value = $newTimestamp
$summary.Range("B6").Value = $newTimestamp

# --- 2) Per-profile detail sheets: append rows 10 & 11 ---
# Each sheet already ends with a pair of rows (8 = "even" style, 9 = "odd"
# style) that log a check. We clone that same two-row pattern for the new
# 22:01 check, then set the new timestamp (and, where relevant, the
# "Szczegóły nowych" listing id string) on the clones.

$detailSheets = @{
    "wszystkie-lublin" = @{ I10 = $null;                                  I11 = $null }
    "artymiuk"         = @{ I10 = $null;                                  I11 = $null }
    "poqui"            = @{ I10 = "183ger|18KAEc|17NeTz|17vbYq|1951OR";   I11 = "17vbYq|17NeTz|183ger|18KAEc|1951OR" }
    "stylowepokoje"    = @{ I10 = "195dLc|16ZeYm";                        I11 = "195dLc|16ZeYm" }
    "villahome"        = @{ I10 = $null;                                  I11 = $null }
}

foreach ($name in $detailSheets.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $info = $detailSheets[$name]

    # Clone row 8 (even style: s=7/8/8/8/9/10/10/5), incl. hidden col I,
    # into row 10
    $ws.Range("A8:I8").Copy($ws.Range("A10:I10"))
    $ws.Rows.Item(10).RowHeight = 18
    $ws.Range("A10").Value = $newTimestamp
    if ($info.I10) {
        $ws.Range("I10").Value = $info.I10
    }

    # Clone row 9 (odd style: s=7/3/3/3/4/12/12/5), incl. hidden col I,
    # into row 11
    $ws.Range("A9:I9").Copy($ws.Range("A11:I11"))
    $ws.Rows.Item(11).RowHeight = 18
    $ws.Range("A11").Value = $newTimestamp
    if ($info.I11) {
        $ws.Range("I11").Value = $info.I11
    }
}

Write-Output "OLX monitor rows updated for $newTimestamp"
